# cryptos.xlsx — "Updated cryptos list" GitHub Actions refresh.
#
# Updates the Price (D) and Volume(1h) (E) columns for the crypto rows,
# and swaps the Polygon/WrappedEther rows (12 <-> 13) including their
# Coin name (B), Link (C), Price (D) and Volume (E) values.
#
# Price values that look like plain decimals (e.g. "241.60", "0.5296")
# are written with a leading apostrophe so Excel stores them as literal
# text instead of silently normalising them to a number (which would
# drop meaningful trailing zeros, e.g. "241.60" -> 241.6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "29.434.80"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.875.19"
$ws.Range("E3").Value = "  +0.88%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.26%  "

# Row 5 - XRP
$ws.Range("D5").Value = "'0.7139"
$ws.Range("E5").Value = "  +1.53%  "

# Row 6 - BNB
$ws.Range("D6").Value = "'241.60"
$ws.Range("E6").Value = "  +1.36%  "

# Row 7 - USDC (price unchanged, volume only)
$ws.Range("E7").Value = "  +0.21%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value = "'0.07897"
$ws.Range("E8").Value = "  -0.07%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "'0.3087"
$ws.Range("E9").Value = "  +1.40%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'25.45"
$ws.Range("E10").Value = "  +3.97%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.08253"
$ws.Range("E11").Value = "  +0.89%  "

# Row 12 - now Polygon (was WrappedEther)
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.7244"
$ws.Range("E12").Value = "  +0.13%  "

# Row 13 - now WrappedEther (was Polygon)
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.865.10"
$ws.Range("E13").Value = "  +8.91%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.255"
$ws.Range("E14").Value = "  +0.70%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "'91.13"
$ws.Range("E15").Value = "  +1.64%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "29.411.76"
$ws.Range("E16").Value = "  +0.37%  "

# Row 17 - Uniswap
$ws.Range("D17").Value = "'5.862"
$ws.Range("E17").Value = "  +0.74%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'244.77"
$ws.Range("E18").Value = "  +2.57%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.000007827"
$ws.Range("E19").Value = "  +0.11%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'13.23"
$ws.Range("E20").Value = "  -0.14%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.120.70"
$ws.Range("E21").Value = "  +13.32%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "'8.048"
$ws.Range("E22").Value = "  +6.50%  "

# Row 23 - Dai (price unchanged, volume only)
$ws.Range("E23").Value = "  +0.09%  "

# Row 24 - BinanceUSD (price unchanged, volume only)
$ws.Range("E24").Value = "  +0.30%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1602"
$ws.Range("E25").Value = "  +11.85%  "

# Row 26 - Monero
$ws.Range("D26").Value = "'162.65"
$ws.Range("E26").Value = "  +0.35%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "'9.005"
$ws.Range("E27").Value = "  +1.12%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'18.29"
$ws.Range("E28").Value = "  +0.86%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'1.357"
$ws.Range("E29").Value = "  -1.95%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "'1.498"
$ws.Range("E30").Value = "  +1.50%  "

# Row 31 - Filecoin (price unchanged, volume only)
$ws.Range("E31").Value = "  +1.67%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "'4.099"
$ws.Range("E32").Value = "  +1.03%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "'0.05194"
$ws.Range("E33").Value = "  +0.11%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "'1.933"
$ws.Range("E34").Value = "  +0.90%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "'1.192"
$ws.Range("E35").Value = "  +1.48%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "'0.7217"
$ws.Range("E36").Value = "  +1.40%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "'2.675"
$ws.Range("E37").Value = "  -0.01%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "'0.01863"
$ws.Range("E38").Value = "  +0.59%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "'2.689"
$ws.Range("E39").Value = "  -0.02%  "

# Row 40 - Maker
$ws.Range("D40").Value = "1.179.96"
$ws.Range("E40").Value = "  +1.83%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.9083"
$ws.Range("E41").Value = "  -1.41%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "'6.127"
$ws.Range("E42").Value = "  +3.19%  "

# Row 43 - Aave
$ws.Range("D43").Value = "'72.53"
$ws.Range("E43").Value = "  +2.01%  "

# Row 44 - PaxDollar (price unchanged, volume only)
$ws.Range("E44").Value = "  +0.25%  "

# Row 45 - Quant
$ws.Range("D45").Value = "'102.07"
$ws.Range("E45").Value = "  +1.06%  "

# Row 46 - Mantle (volume unchanged, price only)
$ws.Range("D46").Value = "'0.5296"

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "2.012.24"
$ws.Range("E47").Value = "  +11.20%  "

# Row 48 - RenderToken
$ws.Range("D48").Value = "'1.792"
$ws.Range("E48").Value = "  +2.03%  "

# Row 49 - SynthetixNetwork
$ws.Range("D49").Value = "'2.904"
$ws.Range("E49").Value = "  +6.45%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'9.289"
$ws.Range("E50").Value = "  +1.03%  "

# Row 51 - TheSandbox
$ws.Range("D51").Value = "'0.4292"
$ws.Range("E51").Value = "  +0.76%  "
